# Fruta / hortaliza, semanal
# Insert two new weekly-price rows at the top of the data block (row 81),
# pushing the existing rows 81-95 down to 83-97.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("81:82").Insert()

# New row 81
$ws.Range("A81").Value = 4
$ws.Range("B81").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C81").Value = 'Los Lagos'
$ws.Range("D81").Value = 44505
$ws.Range("E81").Value = 10
$ws.Range("F81").Value = 'Fruta'
$ws.Range("G81").Value = 100108
$ws.Range("H81").Value = 'Tropicales y subtropicales'
$ws.Range("I81").Value = 100108002
$ws.Range("J81").Value = 'Mango'
$ws.Range("K81").Value = 'Sin especificar'
$ws.Range("L81").Value = 'Primera'
$ws.Range("M81").Value = 200
$ws.Range("N81").Value = 8000
$ws.Range("O81").Value = 8500
$ws.Range("P81").Value = 8250
$ws.Range("Q81").Value = '$/bandeja 4 kilos'
$ws.Range("R81").Value = 'Perú'
$ws.Range("S81").Value = 2062
$ws.Range("T81").Value = 4

# New row 82
$ws.Range("A82").Value = 4
$ws.Range("B82").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C82").Value = 'Los Lagos'
$ws.Range("D82").Value = 44505
$ws.Range("E82").Value = 10
$ws.Range("F82").Value = 'Fruta'
$ws.Range("G82").Value = 100108
$ws.Range("H82").Value = 'Tropicales y subtropicales'
$ws.Range("I82").Value = 100108002
$ws.Range("J82").Value = 'Mango'
$ws.Range("K82").Value = 'Sin especificar'
$ws.Range("L82").Value = 'Segunda'
$ws.Range("M82").Value = 150
$ws.Range("N82").Value = 6000
$ws.Range("O82").Value = 6000
$ws.Range("P82").Value = 6000
$ws.Range("Q82").Value = '$/bandeja 4 kilos'
$ws.Range("R82").Value = 'Perú'
$ws.Range("S82").Value = 1500
$ws.Range("T82").Value = 4
